# Update the "Login" sheet test credentials with the freshly generated
# manager account (username/password) used by the automated test suite,
# and restore the workbook view state so the Login tab is the active /
# selected sheet (matching the "newCustomer" sheet's selection reset).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Login")
$ws2 = $wb.Worksheets.Item("newCustomer")

# New manager login test data (old: mngr511285 / YpepUne)
$ws1.Range("A2").Value = "mngr516939"
$ws1.Range("B2").Value = "dyhAmUd"

# Reset selection on the "newCustomer" sheet and make it inactive
$ws2.Activate()
$ws2.Range("H9").Select()

# Make "Login" the active sheet/tab with its own selection restored
$ws1.Activate()
$ws1.Range("C12").Select()
